# Update countries & provincias Spain
# - Refresh "Estados Unidos" (row 4) totals
# - Insert a new "Venezuela" entry right after "Guatemala" (row 117),
#   pushing Kenia/Mayotte/Vietnam down one row, and drop the old
#   Venezuela row (now redundant, further down the list) so the total
#   row count is unchanged.
# - Refresh the "Datos actualizados..." timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Estados Unidos (row 4) totals ---
$ws.Range("B4").Value = 792759
$ws.Range("C4").Value = 28123
$ws.Range("D4").Value = 72389
$ws.Range("E4").Value = 677856
$ws.Range("F4").Value = 13951
$ws.Range("G4").Value = 1939
$ws.Range("H4").Value = 42514

# --- Insert Venezuela right after Guatemala (row 117) ---
# This shifts Kenia/Mayotte/Vietnam/Venezuela(old)/Tanzania/... down by one row.
$ws.Rows("118:118").Insert()

$ws.Range("A118").Value = "Venezuela"
$ws.Range("B118").Value = 285
$ws.Range("C118").Value = 29
$ws.Range("D118").Value = 117
$ws.Range("E118").Value = 158
$ws.Range("F118").Value = 4
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 10

# The old Venezuela row has now shifted down to row 122 (was row 121).
# Remove it so the overall row count / dimension stays the same.
$ws.Rows("122:122").Delete()

# --- Update the "updated at" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 02:22"
